# Backup-before-publication edit:
# Add a new worksheet "all" (placed first) that concatenates the data
# from the four existing sheets (G3_effective, G3_ineffective,
# notG3_effective, SHH_effective) and adds an explicit "index" header
# for the leading (pandas) index column.

$wb = $excel.ActiveWorkbook

# Keep references to the four existing sheets before inserting the new one.
$sheetG3eff    = $wb.Worksheets.Item("G3_effective")
$sheetG3ineff  = $wb.Worksheets.Item("G3_ineffective")
$sheetNotG3eff = $wb.Worksheets.Item("notG3_effective")
$sheetSHHeff   = $wb.Worksheets.Item("SHH_effective")

# Create the new sheet "all" in front of everything else.
$allSheet = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$allSheet.Name = "all"

# ---- Header row -------------------------------------------------------
$headers = New-Object 'object[,]' 1,11
$headers[0,0]  = "index"
$headers[0,1]  = "Name"
$headers[0,2]  = "G3"
$headers[0,3]  = "mean_rank(G3)"
$headers[0,4]  = "G4"
$headers[0,5]  = "mean_rank(G4)"
$headers[0,6]  = "SHH"
$headers[0,7]  = "mean_rank(SHH)"
$headers[0,8]  = "SHH+p53"
$headers[0,9]  = "mean_rank(SHH+p53)"
$headers[0,10] = "SBI"
$allSheet.Range("A1:K1").Value = $headers
$allSheet.Range("L1").Value = "sbi_id"

# ---- Data rows ----------------------------------------------------------
# Row order: all rows from G3_effective, then G3_ineffective, then
# notG3_effective, then SHH_effective (this matches how the sheets were
# originally concatenated into the combined "all" sheet).
$data = New-Object 'object[,]' 10,12

# G3_effective (old sheet1) rows
$data[0,0]=4;   $data[0,1]="gsk1070916"; $data[0,2]=1; $data[0,3]=2.1;               $data[0,4]=1; $data[0,5]=2.75;   $data[0,6]=1; $data[0,7]=2;                  $data[0,8]=1; $data[0,9]=2;                  $data[0,10]="['SBI-0654453.P001']"; $data[0,11]="SBI-0654453.P001"
$data[1,0]=11;  $data[1,1]="bx-912";     $data[1,2]=1; $data[1,3]=10.4;              $data[1,4]=1; $data[1,5]=8.75;   $data[1,6]=1; $data[1,7]=116;                $data[1,8]=0; $data[1,9]=223;                $data[1,10]="['SBI-0645949.P001']"; $data[1,11]="SBI-0645949.P001"
$data[2,0]=44;  $data[2,1]="abt-737";    $data[2,2]=1; $data[2,3]=32.049999999999997;$data[2,4]=1; $data[2,5]=32.75;  $data[2,6]=1; $data[2,7]=33.416666666666657; $data[2,8]=1; $data[2,9]=33.833333333333343; $data[2,10]="[]";                   $data[2,11]="NaN"
$data[3,0]=108; $data[3,1]="linsitinib"; $data[3,2]=1; $data[3,3]=84.35;             $data[3,4]=1; $data[3,5]=88.75;  $data[3,6]=1; $data[3,7]=87.25;              $data[3,8]=1; $data[3,9]=86.5;               $data[3,10]="['SBI-0646932.P001']"; $data[3,11]="SBI-0646932.P001"

# G3_ineffective (old sheet2) rows
$data[4,0]=135; $data[4,1]="etoposide";  $data[4,2]=0; $data[4,3]=156.94999999999999;$data[4,4]=0; $data[4,5]=167.625;$data[4,6]=0; $data[4,7]=163.5;              $data[4,8]=0; $data[4,9]=184.66666666666671;$data[4,10]="['SBI-0634393.P001', 'SBI-0634371.P001', 'SBI-0051910.P001', 'SBI-0050405.P002', 'SBI-0634371.P002', 'SBI-0051910.P002']"; $data[4,11]="SBI-0634393.P001"
$data[5,0]=149; $data[5,1]="topotecan";  $data[5,2]=0; $data[5,3]=228.05;            $data[5,4]=0; $data[5,5]=217.25; $data[5,6]=0; $data[5,7]=214.66666666666671; $data[5,8]=0; $data[5,9]=189;               $data[5,10]="['SBI-0055592.P003', 'SBI-0055592.P005', 'SBI-0055592.P002']"; $data[5,11]="SBI-0055592.P003"

# notG3_effective (old sheet3) rows
$data[6,0]=12;  $data[6,1]="vx-702";     $data[6,2]=1; $data[6,3]=7.9;               $data[6,4]=1; $data[6,5]=8.5;    $data[6,6]=1; $data[6,7]=7.666666666666667;  $data[6,8]=1; $data[6,9]=7;                  $data[6,10]="['SBI-0654264.P001']"; $data[6,11]="SBI-0654264.P001"
$data[7,0]=64;  $data[7,1]="alisertib";  $data[7,2]=1; $data[7,3]=34.299999999999997;$data[7,4]=1; $data[7,5]=48.25;  $data[7,6]=1; $data[7,7]=55.25;              $data[7,8]=1; $data[7,9]=42.5;               $data[7,10]="['SBI-0646927.P001']"; $data[7,11]="SBI-0646927.P001"

# SHH_effective (old sheet4) rows
$data[8,0]=35;  $data[8,1]="olaparib";   $data[8,2]=0; $data[8,3]=155.9;             $data[8,4]=0; $data[8,5]=113.75; $data[8,6]=1; $data[8,7]=23;                 $data[8,8]=1; $data[8,9]=22.833333333333329;$data[8,10]="[]"; $data[8,11]="NaN"
$data[9,0]=40;  $data[9,1]="rucaparib";  $data[9,2]=0; $data[9,3]=113.4;             $data[9,4]=1; $data[9,5]=21.75;  $data[9,6]=1; $data[9,7]=15.75;              $data[9,8]=1; $data[9,9]=16.166666666666671;$data[9,10]="[]"; $data[9,11]="NaN"

$allSheet.Range("A2:L11").Value = $data

# ---- Formatting ---------------------------------------------------------
# Column A (index) and the header row use the same bold / centered /
# thin-bordered style that the "Name" ... "sbi_id" headers and the other
# sheets' leading index column already use.
$headerStyleRange = $allSheet.Range("A1:L1")
$headerStyleRange.Font.Bold = $true
$headerStyleRange.HorizontalAlignment = -4108   # xlCenter
$headerStyleRange.VerticalAlignment = -4160     # xlTop
$headerStyleRange.Borders.LineStyle = 1         # xlContinuous
$headerStyleRange.Borders.Weight = 2            # xlThin

$indexColRange = $allSheet.Range("A2:A11")
$indexColRange.Font.Bold = $true
$indexColRange.HorizontalAlignment = -4108
$indexColRange.VerticalAlignment = -4160
$indexColRange.Borders.LineStyle = 1
$indexColRange.Borders.Weight = 2

# Select C11 to mirror the saved selection state for this sheet.
$allSheet.Range("C11").Select() | Out-Null
